$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $r = $d.Content
    $ok = $r.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $ok) {
        Write-Output "WARNING: not found: $oldText"
    }
}

# --- Title / byline / contact line -----------------------------------
Replace-Text "Unveiling the Convergence of Science and Art" "Beyond the Classroom: The Value of Arts Education"
Replace-Text "Ethan James" "Jessica Albright"
Replace-Text "ethan" "jessica"
Replace-Text "james@emailcentral" "albright@schoolmail"
Replace-Text "xyz" "com"

# --- Body paragraph (essay) -------------------------------------------
Replace-Text "The intersection of science and art is a realm where boundaries blur, and creativity and knowledge converge" "The arts, in their myriad forms, have long been recognized for their ability to enrich our lives and expand our horizons"
Replace-Text " It is a space where artistic expression and scientific inquiry dance in harmonious unison, each enriching the other" " As educators, it is our duty to cultivate a comprehensive educational experience that values the arts alongside the traditional subjects"
Replace-Text " This essay delves into the fascinating landscape of this convergence, exploring how science inspires art and how art, in turn, informs science" " This essay explores the invaluable role of arts education in shaping well-rounded, creative, and expressive individuals"
Replace-Text "In the tapestry of this convergence, art finds its muse in the wonders of the natural world and the intricacies of scientific phenomena" "In a world increasingly dominated by technology and quantitative reasoning, the arts offer a refuge for the imaginative mind"
Replace-Text " Artists draw inspiration from the colors of the aurora borealis, the patterns of snowflakes, and the intricate geometry of DNA" " Through engagement with visual arts, music, theater, and dance, students can cultivate their unique perspectives and find creative outlets for self-expression"
Replace-Text " These natural marvels become the raw materials of artistic creation, transformed into paintings, sculptures, and installations that capture the essence of scientific concepts" " The arts provide a safe haven for exploration, where students can experiment with different forms of expression, allowing them to develop their own artistic voices"
Replace-Text "Conversely, science also finds a muse in art" "Furthermore, the arts foster essential critical thinking and problem-solving skills"
Replace-Text " The creative process inherent in art encourages scientists to think outside the conventional boundaries of their disciplines" " In analyzing and interpreting works of art, students develop the ability to think critically and engage in meaningful discourse"
Replace-Text " Artists' unique perspectives and methodologies prompt scientists to approach problems from novel angles, leading to breakthroughs and innovations" " They learn to appreciate different perspectives and understand the nuances of communication"
Replace-Text " This reciprocal relationship between science and art fosters a fertile ground for interdisciplinary exploration and discovery" " By actively participating in the creative process, students learn to embrace challenges, think outside the box, and find innovative solutions to problems"

# Remove the trailing sentences that were deleted from the essay paragraph
# (everything after "...find innovative solutions to problems.")
$target = " By actively participating in the creative process, students learn to embrace challenges, think outside the box, and find innovative solutions to problems"
$r = $d.Content
$found = $r.Find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $periodEnd = $r.End + 1
    $essayPara = $d.Paragraphs.Item(5)
    $delEnd = $essayPara.Range.End - 1
    if ($delEnd -gt $periodEnd) {
        $delRange = $d.Range($periodEnd, $delEnd)
        $delRange.Delete()
    }
}

# --- Summary paragraph --------------------------------------------------
Replace-Text "The convergence of science and art is a dynamic and ever-evolving realm where creativity and knowledge intertwine" "In essence, arts education provides students with the tools and skills necessary to navigate the complexities of an ever-changing world"
Replace-Text " Science inspires art, providing a wealth of ideas and wonders that fuel artistic expression" " It cultivates creativity, critical thinking, and empathy, while fostering a sense of personal and cultural identity"
Replace-Text " Simultaneously, art informs science, challenging conventional thinking and stimulating innovative approaches to problem-solving" " By recognizing the importance of the arts in education, we empower students to become well-rounded individuals who are equipped to make meaningful contributions to society"

# Remove the trailing sentence deleted from the summary paragraph
$target2 = " By recognizing the importance of the arts in education, we empower students to become well-rounded individuals who are equipped to make meaningful contributions to society"
$r2 = $d.Content
$found2 = $r2.Find.Execute($target2, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $periodEnd2 = $r2.End + 1
    $summaryPara = $d.Paragraphs.Last
    $delEnd2 = $summaryPara.Range.End - 1
    if ($delEnd2 -gt $periodEnd2) {
        $delRange2 = $d.Range($periodEnd2, $delEnd2)
        $delRange2.Delete()
    }
}

# --- Fix font typo across the whole document ----------------------------
foreach ($p in $d.Paragraphs) {
    $p.Range.Font.Name = "Times New Roman"
}

# --- Append a new trailing empty paragraph -------------------------------
$d.Paragraphs.Last.Range.InsertParagraphAfter()

Write-Output "done"
